# The deck's single slide master (ppt/theme/theme1.xml) currently uses the
# "Integral" theme's "Red Violet" colour scheme. Re-colour it with the
# stock "Office" palette (the scheme that, before this edit, only the
# notes master's otherwise-unused theme part carried) via the Design >
# Colors object model, i.e. ThemeColorScheme.Colors(n).RGB, which is the
# supported COM surface for rewriting a theme's <a:clrScheme> entries.
#
# OLE_COLOR/RGB() values are byte-order BGR (0x00BBGGRR), so convert each
# target hex swatch (RRGGBB, as authored in DrawingML) accordingly.
function ConvertTo-OleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (theme colour slot, target "Office" hex value)
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Colors($i).RGB = ConvertTo-OleColor $officeColors[$i - 1]
}
